$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 495
$ws.Range("I61").Value = 495
$ws.Range("K61").Value = 1485
$ws.Range("M61").Value = -1313

$ws.Range("H62").Value = 12386.375
$ws.Range("I62").Value = 13182
$ws.Range("J62").Value = 9999.5
$ws.Range("K62").Value = 13182
$ws.Range("L62").Value = 9999.5
$ws.Range("M62").Value = -12558
$ws.Range("N62").Value = -11247.5

$ws.Range("H65").Value = 12386.375
$ws.Range("I65").Value = 13182
$ws.Range("J65").Value = 9999.5
$ws.Range("K65").Value = 65910
$ws.Range("L65").Value = 49997.5
$ws.Range("M65").Value = -62790
$ws.Range("N65").Value = -56237.5

$ws.Range("H98").Value = 1341.7391
$ws.Range("J98").Value = 987.5
$ws.Range("L98").Value = 987.5
$ws.Range("N98").Value = -3983.5

$ws.Range("H101").Value = 55558136
$ws.Range("I101").Value = 71429030
$ws.Range("J101").Value = 9999
$ws.Range("K101").Value = 214287090
$ws.Range("L101").Value = 29997
$ws.Range("M101").Value = -214285468
$ws.Range("N101").Value = -33241

$ws.Range("H113").Value = 6440.7144
$ws.Range("I113").Value = 5150
$ws.Range("J113").Value = 6957
$ws.Range("K113").Value = 5150
$ws.Range("L113").Value = 6957
$ws.Range("M113").Value = -1896
$ws.Range("N113").Value = -13465

$ws.Range("H122").Value = 1341.7391
$ws.Range("J122").Value = 987.5
$ws.Range("L122").Value = 2962.5
$ws.Range("N122").Value = -7862.5

$ws.Range("H133").Value = 50709
$ws.Range("I133").Value = 50709
$ws.Range("K133").Value = 50709
$ws.Range("M133").Value = -45649

$ws.Range("H138").Value = 2497.4443
$ws.Range("J138").Value = 4472.2964
$ws.Range("L138").Value = 13416.8892
$ws.Range("N138").Value = -23696.8892

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 883.1667
$ws.Range("J4").Value = 383.66666
$ws.Range("L4").Value = 383.66666
$ws.Range("N4").Value = -615.66666

$ws.Range("H32").Value = 3725.375
$ws.Range("I32").Value = 2303
$ws.Range("K32").Value = 2303
$ws.Range("M32").Value = -2016

$ws.Range("H122").Value = 1491358.4
$ws.Range("I122").Value = 3547.8
$ws.Range("K122").Value = 10643.4
$ws.Range("M122").Value = -8193.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3254292
$ws.Range("I94").Value = 4546398
$ws.Range("J94").Value = 24027.5
$ws.Range("K94").Value = 4546398
$ws.Range("L94").Value = 24027.5
$ws.Range("M94").Value = -4545947
$ws.Range("N94").Value = -24929.5

$ws.Range("H107").Value = 2235554
$ws.Range("I107").Value = 2977997.2
$ws.Range("J107").Value = 8224.625
$ws.Range("K107").Value = 2977997.2
$ws.Range("L107").Value = 8224.625
$ws.Range("M107").Value = -2976077.2
$ws.Range("N107").Value = -12064.625

$ws.Range("H115").Value = 48000
$ws.Range("J115").Value = 48000
$ws.Range("L115").Value = 48000
$ws.Range("N115").Value = -51134

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 299.23077
$ws.Range("I7").Value = 197.625
$ws.Range("K7").Value = 197.625
$ws.Range("M7").Value = -84.625

$ws.Range("H31").Value = 13561.826
$ws.Range("I31").Value = 1927.6923
$ws.Range("K31").Value = 1927.6923
$ws.Range("M31").Value = -1632.6923

$ws.Range("H34").Value = 13561.826
$ws.Range("I34").Value = 1927.6923
$ws.Range("K34").Value = 1927.6923
$ws.Range("M34").Value = -1725.6923

$ws.Range("H86").Value = 8776.450000000001
$ws.Range("I86").Value = 7963.875
$ws.Range("J86").Value = 9318.166999999999
$ws.Range("K86").Value = 7963.875
$ws.Range("L86").Value = 9318.166999999999
$ws.Range("M86").Value = -6840.875
$ws.Range("N86").Value = -11564.167

$ws.Range("H89").Value = 8776.450000000001
$ws.Range("I89").Value = 7963.875
$ws.Range("J89").Value = 9318.166999999999
$ws.Range("K89").Value = 39819.375
$ws.Range("L89").Value = 46590.835
$ws.Range("M89").Value = -34203.375
$ws.Range("N89").Value = -57822.835

$ws.Range("H99").Value = 4527.222
$ws.Range("I99").Value = 5500
$ws.Range("J99").Value = 4249.2856
$ws.Range("K99").Value = 5500
$ws.Range("L99").Value = 4249.2856
$ws.Range("M99").Value = -4002
$ws.Range("N99").Value = -7245.2856

$ws.Range("H126").Value = 4527.222
$ws.Range("I126").Value = 5500
$ws.Range("J126").Value = 4249.2856
$ws.Range("K126").Value = 16500
$ws.Range("L126").Value = 12747.8568
$ws.Range("M126").Value = -14030
$ws.Range("N126").Value = -17687.8568

$ws.Range("H132").Value = 120493
$ws.Range("I132").Value = 79269.234
$ws.Range("K132").Value = 237807.702
$ws.Range("M132").Value = -235277.702

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4472633
$ws.Range("I4").Value = 8124930
$ws.Range("J4").Value = 89877.13
$ws.Range("K4").Value = 24374790
$ws.Range("L4").Value = 269631.39
$ws.Range("M4").Value = -24374678
$ws.Range("N4").Value = -269855.39

$ws.Range("H12").Value = 55779.438
$ws.Range("J12").Value = 324.16666
$ws.Range("L12").Value = 972.4999799999999
$ws.Range("N12").Value = -1318.49998

$ws.Range("H37").Value = 63671.43
$ws.Range("J37").Value = 63671.43
$ws.Range("L37").Value = 191014.29
$ws.Range("N37").Value = -191238.29

$ws.Range("H69").Value = 4257.5
$ws.Range("I69").Value = 1012
$ws.Range("J69").Value = 4906.6
$ws.Range("K69").Value = 3036
$ws.Range("L69").Value = 14719.8
$ws.Range("M69").Value = -2225
$ws.Range("N69").Value = -16341.8

$ws.Range("H72").Value = 4257.5
$ws.Range("I72").Value = 1012
$ws.Range("J72").Value = 4906.6
$ws.Range("K72").Value = 9108
$ws.Range("L72").Value = 44159.4
$ws.Range("M72").Value = -5052
$ws.Range("N72").Value = -52271.4

$ws.Range("H92").Value = 875.63635
$ws.Range("J92").Value = 967
$ws.Range("L92").Value = 2901
$ws.Range("N92").Value = -5397

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 6499.9443
$ws.Range("I11").Value = 4799.933
$ws.Range("K11").Value = 4799.933
$ws.Range("M11").Value = -4660.933

$ws.Range("H13").Value = 952.5
$ws.Range("I13").Value = 952.5
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 952.5
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -813.5
$ws.Range("N13").ClearContents()

$ws.Range("H97").Value = 1036840.8
$ws.Range("I97").Value = 1985787.1
$ws.Range("J97").Value = 1626.5454
$ws.Range("K97").Value = 1985787.1
$ws.Range("L97").Value = 1626.5454
$ws.Range("M97").Value = -1985291.1
$ws.Range("N97").Value = -2618.5454

$ws.Range("H102").Value = 5399638
$ws.Range("I102").Value = 6537696
$ws.Range("K102").Value = 6537696
$ws.Range("M102").Value = -6536074

$ws.Range("H122").Value = 1485693.5
$ws.Range("I122").Value = 1485693.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4457080.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4454630.5
$ws.Range("N122").ClearContents()

$ws.Range("H123").Value = 50000
$ws.Range("J123").Value = 50000
$ws.Range("L123").Value = 50000
$ws.Range("N123").Value = -54900

$ws.Range("H126").Value = 4444186
$ws.Range("I126").Value = 2069266.6
$ws.Range("J126").Value = 11908219
$ws.Range("K126").Value = 6207799.800000001
$ws.Range("L126").Value = 35724657
$ws.Range("M126").Value = -6205329.800000001
$ws.Range("N126").Value = -35729597

$ws.Range("H132").Value = 3384.9312
$ws.Range("J132").Value = 4733.1
$ws.Range("L132").Value = 14199.3
$ws.Range("N132").Value = -19259.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 12005620
$ws.Range("I2").Value = 15003275
$ws.Range("J2").Value = 15000
$ws.Range("K2").Value = 15003275
$ws.Range("L2").Value = 15000
$ws.Range("M2").Value = -15003163
$ws.Range("N2").Value = -15224

$ws.Range("H22").Value = 128798.14
$ws.Range("I22").Value = 296929.34
$ws.Range("J22").Value = 2699.75
$ws.Range("K22").Value = 296929.34
$ws.Range("L22").Value = 2699.75
$ws.Range("M22").Value = -296634.34
$ws.Range("N22").Value = -3289.75

$ws.Range("H27").Value = 128798.14
$ws.Range("I27").Value = 296929.34
$ws.Range("J27").Value = 2699.75
$ws.Range("K27").Value = 296929.34
$ws.Range("L27").Value = 2699.75
$ws.Range("M27").Value = -296822.34
$ws.Range("N27").Value = -2913.75

$ws.Range("H40").Value = 5573.222
$ws.Range("I40").Value = 3782
$ws.Range("K40").Value = 3782
$ws.Range("M40").Value = -3646

$ws.Range("H46").Value = 9908.727999999999
$ws.Range("I46").Value = 7928.4287
$ws.Range("K46").Value = 7928.4287
$ws.Range("M46").Value = -7740.4287

$ws.Range("H55").Value = 2378.611
$ws.Range("I55").Value = 1917.5834
$ws.Range("K55").Value = 1917.5834
$ws.Range("M55").Value = -1744.5834

$ws.Range("H93").Value = 11906410
$ws.Range("I93").Value = 14494440
$ws.Range("K93").Value = 14494440
$ws.Range("M93").Value = -14493192

$ws.Range("H122").Value = 7134.6924
$ws.Range("J122").Value = 7926.7
$ws.Range("L122").Value = 23780.1
$ws.Range("N122").Value = -28680.1

$ws.Range("H136").Value = 128806.625
$ws.Range("I136").Value = 169825.58
$ws.Range("K136").Value = 509476.74
$ws.Range("M136").Value = -506926.74

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2653.12
$ws.Range("I126").Value = 2286.5
$ws.Range("K126").Value = 6859.5
$ws.Range("M126").Value = -4389.5

$ws.Range("H132").Value = 38885796
$ws.Range("I132").Value = 41672576
$ws.Range("K132").Value = 125017728
$ws.Range("M132").Value = -125015198

$ws.Range("H136").Value = 5342.857
$ws.Range("J136").Value = 6854.5454
$ws.Range("L136").Value = 20563.6362
$ws.Range("N136").Value = -25663.6362
